# Auto update Excel log
# Appends newly-detected sensor events to the ALERTS log and to the
# mmWave raw-event log.
#
# Note: plain-looking "date" strings like "2026-01-30" get auto-converted
# to Excel date serials by Range.Value assignment. The source log stores
# these as literal text, so we briefly force a text number format, assign
# the value, then clear the formatting again so the cell ends up as plain
# text with no residual style (matching every other cell in these sheets).

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($Sheet, $Row, $Date, $Timestamp, $Hour, $Location, $Value, $Status)

    $dateCell = $Sheet.Cells.Item($Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $Date
    $dateCell.ClearFormats()

    $Sheet.Cells.Item($Row, 2).Value = $Timestamp
    $Sheet.Cells.Item($Row, 3).Value = $Hour
    $Sheet.Cells.Item($Row, 4).Value = $Location
    $Sheet.Cells.Item($Row, 5).Value = $Value
    $Sheet.Cells.Item($Row, 6).Value = $Status
}

# --- ALERTS sheet: new fall-detection critical emergency row ---
$alerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $alerts 12 "2026-01-30" "14:19:45" "14:00" "Living Room" "FALL_DETECTED" "CRITICAL EMERGENCY"

# --- mmWave sheet: raw presence/fall/presence sequence around the fall ---
$mmwave = $wb.Worksheets.Item("mmWave")
Add-LogRow $mmwave 65 "2026-01-30" "14:19:25" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmwave 66 "2026-01-30" "14:19:45" "14:00" "Living Room" "FALL_DETECTED" "CRITICAL EMERGENCY"
Add-LogRow $mmwave 67 "2026-01-30" "14:20:23" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
